# Updated symbol list on Wed Dec 21 23:43:30 UTC 2022 with GitHub Actions
#
# The "Price" column (D) holds numeric-looking quotes as TEXT (inline strings
# in the original file), so each one must be written back as text (via
# NumberFormat "@") to avoid Excel silently re-casting it to a float and
# mangling trailing/leading zeros (e.g. "0.05696" -> "0.05700").
# The "Volume(1h)" column (E) cells are plain text labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

Set-TextCell "D2"  "246.67"      # was 246.77
Set-TextCell "D3"  "22.66"       # was 22.68
Set-TextCell "D4"  "5.254"       # was 5.253
Set-TextCell "D5"  "0.05700"     # was 0.05696
Set-TextCell "D6"  "3.415"       # was 3.416
Set-TextCell "D7"  "6.293"       # was 6.291
Set-TextCell "D8"  "0.8102"      # was 0.8097
Set-TextCell "D9"  "0.8714"      # was 0.8689
Set-TextCell "D10" "0.1421"      # was 0.1426
Set-TextCell "D11" "0.07345"     # was 0.07354

Set-TextCell "D13" "0.03076"     # was 0.03078
Set-TextCell "D14" "0.09388"     # was 0.09391
Set-TextCell "D15" "3.877"       # was 3.869
Set-TextCell "D16" "0.001569"    # was 0.001575
Set-TextCell "D17" "0.04797"     # was 0.04786

Set-TextCell "D18" "0.0005839"   # was 0.0005851
$ws.Range("E18").Value = "17OneONEWorstin24h"   # was 17OneONE

Set-TextCell "D19" "0.006083"    # was 0.006081
Set-TextCell "D20" "0.005024"    # was 0.005025
Set-TextCell "D21" "0.0009962"   # was 0.0009969

Set-TextCell "D23" "3.688"       # was 3.690
Set-TextCell "D24" "2.195"       # was 2.196
Set-TextCell "D25" "0.3259"      # was 0.3257

Set-TextCell "D40" "0.03926"     # was 0.03928
Set-TextCell "D41" "0.006748"    # was 0.006789
Set-TextCell "D42" "0.1067"      # was 0.1068
Set-TextCell "D43" "0.003199"    # was 0.002679
Set-TextCell "D44" "0.007502"    # was 0.007504
Set-TextCell "D45" "0.00005593"  # was 0.00005595

Set-TextCell "D47" "0.6108"      # was 0.6111

Set-TextCell "D48" "0.1904"      # was 0.1928
$ws.Range("E48").Value = "47BOLOBOLO"           # was 47BOLOBOLOWorstin24h
